# Update countries & provincias Spain
# The source data feed reshuffled a handful of country rows (alphabetical /
# ranking reorder touched the shared-string table) and refreshed the day's
# case counters. This reproduces the resulting cell values on the "Pais"
# sheet, plus the "last updated" timestamp in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CountryRow {
    param($Row, $Pais, $CasosTotales, $NuevosCasos, $CasosActivos, $Recuperados, $CasosCriticos, $MuertesHoy, $Muertes)

    $ws.Range("A$Row").Value = $Pais
    $ws.Range("B$Row").Value = $CasosTotales
    $ws.Range("C$Row").Value = $NuevosCasos
    $ws.Range("D$Row").Value = $CasosActivos
    $ws.Range("E$Row").Value = $Recuperados
    $ws.Range("F$Row").Value = $CasosCriticos
    $ws.Range("G$Row").Value = $MuertesHoy
    $ws.Range("H$Row").Value = $Muertes
}

# --- Header: last-updated timestamp -------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 20 de Mayo de 2020 a las 01:35"

# --- Row data: País, Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos críticos, Muertes hoy, Muertes --------------------

# Row 4: Estados Unidos (counts refreshed, no reorder)
Set-CountryRow 4 "Estados Unidos" 1570144 19850 363069 1113564 0 1530 93511

# Row 17: Canada (counts refreshed, no reorder)
Set-CountryRow 17 "Canada" 79112 1040 40050 33150 0 70 5912

# Row 50: Panama (counts refreshed, no reorder)
Set-CountryRow 50 "Panama" 9867 261 6194 3392 0 6 281

# Rows 51-52: Chequia / Argentina swap order, counts refreshed
Set-CountryRow 51 "Argentina" 8809 438 2872 5544 0 11 393
Set-CountryRow 52 "Chequia" 8647 61 5726 2619 0 5 302

# Rows 61-64: Finlandia / Moldavia / Nigeria reorder, counts refreshed
Set-CountryRow 61 "Nigeria" 6401 226 1734 4475 0 1 192
Set-CountryRow 62 "Finlandia" 6399 19 5000 1098 0 1 301
Set-CountryRow 63 "Moldavia" 6340 202 2508 3611 0 4 221
Set-CountryRow 64 "Ghana" 6096 361 1773 4292 0 2 31

# Rows 75-76: Grecia / Uzbekistan swap order, counts refreshed
Set-CountryRow 75 "Uzbekistan" 2855 64 2338 504 0 0 13
Set-CountryRow 76 "Grecia" 2840 4 1374 1301 0 0 165

# Rows 91-92: Lituania / Republica de Yibuti swap order, counts refreshed
Set-CountryRow 91 "Republica de Yibuti" 1618 100 1033 578 0 0 7
Set-CountryRow 92 "Lituania" 1562 15 1025 477 0 1 60

# Row 119: Uruguay (counts refreshed, no reorder)
Set-CountryRow 119 "Uruguay" 738 1 579 139 0 0 20

# Rows 127-129: Haiti / Jamaica / Sierra Leona reorder, counts refreshed
Set-CountryRow 127 "Sierra Leona" 534 15 167 334 0 0 33
Set-CountryRow 128 "Haiti" 533 0 21 491 0 0 21
Set-CountryRow 129 "Jamaica" 520 0 145 366 0 0 9

# Row 165: Bermudas (counts refreshed, no reorder)
Set-CountryRow 165 "Bermudas" 125 0 78 38 0 0 9
